# Updates cryptos list prices / hourly volume deltas, and swaps in
# MultiversX for FraxShare at the bottom of the table (row 51), plus a
# reordering of VeChain/RenderToken (rows 39-40).
# D/E columns hold text (not numeric) values in the source data -- some
# of the new values parse as plain numbers (e.g. "1.00", "53.50"), so we
# force the cell to Text format before writing, otherwise Excel's COM
# layer auto-coerces the string to a number and the trailing zeros /
# exact literal text would be lost.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.265.18'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.270.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.84%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.58'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.01'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.566'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.20'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.21'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.614.89'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.265.87'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.810'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.139.97'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +13.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0916'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.02'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.36'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.94'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.89'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '41.34'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +11.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.56'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.63'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.53'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.51'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.09%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.87%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.91%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.93'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.52%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0311'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.21'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.45'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -9.35%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.95'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +11.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.762.44'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.27%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '70.02'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.17'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '95.28'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.73'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.50'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.15%  '
